# Fruta / hortaliza, semanal
# Insert a new weekly record for "Terminal La Palmera de La Serena - Perejil"
# as row 77, pushing the former rows 77 and 78 down to 78 and 79.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 77 (shifts old rows 77-78 down to 78-79,
# copying formatting -- including the date style on column D -- from the row above).
$ws.Rows.Item(77).Insert()

# Populate the newly inserted row 77 with the new weekly reading.
$ws.Range("A77").Value = 8
$ws.Range("B77").Value = "Terminal La Palmera de La Serena"
$ws.Range("C77").Value = "Coquimbo"
$ws.Range("D77").Value2 = 44448
$ws.Range("E77").Value = 4
$ws.Range("F77").Value = 100112044
$ws.Range("G77").Value = "Perejil"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 3200
$ws.Range("K77").Value = 2000
$ws.Range("L77").Value = 2500
$ws.Range("M77").Value = 2250
$ws.Range("N77").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O77").Value = "Provincia del Elquí"
$ws.Range("P77").Value = 1500
$ws.Range("Q77").Value = 1.5
$ws.Range("R77").Value = "Hortaliza"
